$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlCenter = -4108

function Add-DataRow {
    param($RowNum, $KeyName)

    # Column A: key name, Text number format, left/default aligned (style matches existing "3")
    $ws.Cells.Item($RowNum, 1).Value = $KeyName
    $ws.Cells.Item($RowNum, 1).NumberFormat = "@"

    # Column B: email, hyperlinked, centered Text (style matches existing "2")
    $ws.Cells.Item($RowNum, 2).Value = "abv@abv.bg"
    $ws.Hyperlinks.Add($ws.Cells.Item($RowNum, 2), "mailto:abv@abv.bg")
    $ws.Range("B10").Copy()
    $ws.Cells.Item($RowNum, 2).PasteSpecial($xlPasteFormats)

    # Column C: String.Empty, centered Text (style matches existing "1")
    $ws.Cells.Item($RowNum, 3).Value = "String.Empty"
    $ws.Range("C10").Copy()
    $ws.Cells.Item($RowNum, 3).PasteSpecial($xlPasteFormats)

    # Column D: numeric 123, centered Text (style matches existing "1")
    $ws.Cells.Item($RowNum, 4).Value = 123
    $ws.Range("D10").Copy()
    $ws.Cells.Item($RowNum, 4).PasteSpecial($xlPasteFormats)

    # Column E: String.Empty, centered Text (style matches existing "1")
    $ws.Cells.Item($RowNum, 5).Value = "String.Empty"
    $ws.Range("E10").Copy()
    $ws.Cells.Item($RowNum, 5).PasteSpecial($xlPasteFormats)
}

Add-DataRow 12 "Login"
Add-DataRow 13 "CreateNewPost"
Add-DataRow 14 "DeletePost"

$ws.Range("F13").Select()

Write-Output "Data driven rows added"
